$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a serial date value that was bumped by one day
# (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04) for every data row (2..387).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 387 }

$ws.Range("C2:C$lastRow").Value = 45203
